# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (column G) values, keyed by row number (row 1 is the header row).
$kValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 1
    10 = 3
    11 = 2
    12 = 1
    13 = 0
    14 = 3
    15 = 0
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 0
    23 = 1
    24 = 3
    25 = 2
    26 = 1
    27 = 0
    28 = 2
    29 = 1
    30 = 2
    31 = 1
    32 = 0
    33 = 2
    34 = 0
    36 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
